# Add a new "localdb" command-group column to the hidden '#system' sheet,
# insert it (alphabetically) into the "target" list, populate its function
# names, and fix up the defined names that pointed at the shifted columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Insert a brand-new column at N; everything from N..AC slides right to O..AD.
$ws.Columns("N").Insert()

# 2) Populate the new "localdb" column with its header + function names.
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# 3) Insert "localdb" into the alphabetical "target" list in column A (row 14,
#    between "jms" and "macro"), pushing the remaining entries down one row.
for ($row = 29; $row -ge 14; $row--) {
    $val = $ws.Cells.Item($row, 1).Value()
    $ws.Cells.Item($row + 1, 1).Value = $val
}
$ws.Cells.Item(14, 1).Value = "localdb"

# 4) Fix up the defined names whose target ranges shifted because of the new
#    column / the extra "target" row.
$wb.Names.Item("macro").RefersTo = "='#system'!`$O`$2:`$O`$4"
$wb.Names.Item("mail").RefersTo = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$21"

# 5) Add the new "localdb" defined name (appended at the end, like the others).
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")

# 6) Restore the originally-active sheet/tab selection (editing the hidden
#    '#system' sheet above made the engine re-derive this otherwise).
$wb.Worksheets.Item("Number_Command_Validation").Activate()
